# Insert a new weekly price record at the top of the data block (row 49),
# pushing all existing records (rows 49-168) down by one row. The last
# existing record (old row 168) ends up in the newly created row 169.
#
# Only the columns that actually vary per record (D = Fecha, J = Volumen,
# K = Precio minimo, L = Precio maximo, M = Precio promedio ponderado,
# O = Origen, P = Precio $/Kg) are shifted; the remaining columns
# (A,B,C,E,F,G,H,I,N,Q,R) are constant for every row in this block, so the
# row-169 clone picks them up automatically when we duplicate row 168.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 49
$lastRow = 168
$newLastRow = $lastRow + 1

# 1) Clone the last row (168) into the brand new row (169) so every column
#    - including the ones that never change - exists (with correct
#    formatting/number formats) before we cascade the shift below.
$ws.Range("A$lastRow`:R$lastRow").Copy($ws.Range("A$newLastRow`:R$newLastRow"))

# 2) Cascade existing records down by one row: for each row from the bottom
#    up to firstRow+1, pull the varying columns from the row above (which
#    still holds its original, un-shifted values at this point).
for ($r = $newLastRow; $r -gt $firstRow; $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value2 = $ws.Cells.Item($src, 4).Value2   # D Fecha
    $ws.Cells.Item($r, 10).Value2 = $ws.Cells.Item($src, 10).Value2 # J Volumen
    $ws.Cells.Item($r, 11).Value2 = $ws.Cells.Item($src, 11).Value2 # K Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $ws.Cells.Item($src, 12).Value2 # L Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $ws.Cells.Item($src, 13).Value2 # M Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value2 = $ws.Cells.Item($src, 15).Value2 # O Origen
    $ws.Cells.Item($r, 16).Value2 = $ws.Cells.Item($src, 16).Value2 # P Precio $/Kg
}

# 3) Write the brand new record into row 49 (Origen stays as-is).
$ws.Cells.Item($firstRow, 4).Value2 = 44498   # D Fecha
$ws.Cells.Item($firstRow, 10).Value2 = 105    # J Volumen
$ws.Cells.Item($firstRow, 11).Value2 = 8000   # K Precio minimo
$ws.Cells.Item($firstRow, 12).Value2 = 8500   # L Precio maximo
$ws.Cells.Item($firstRow, 13).Value2 = 8262   # M Precio promedio ponderado
$ws.Cells.Item($firstRow, 16).Value2 = 138    # P Precio $/Kg
